$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 6874.75
$ws.Range("I82").Value = 1040
$ws.Range("J82").Value = 11042.429
$ws.Range("K82").Value = 3120
$ws.Range("L82").Value = 33127.287
$ws.Range("M82").Value = -2714
$ws.Range("N82").Value = -33939.287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 6874.75
$ws.Range("I85").Value = 1040
$ws.Range("J85").Value = 11042.429
$ws.Range("K85").Value = 3120
$ws.Range("L85").Value = 33127.287
$ws.Range("M85").Value = -1716
$ws.Range("N85").Value = -35935.287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1275.3103
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1275.3103
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 3825.9309
$ws.Range("M112").Value = $null
$ws.Range("N112").Value = -6041.9309

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 827.65
$ws.Range("I129").Value = 355.125
$ws.Range("J129").Value = 868.73914
$ws.Range("K129").Value = 1065.375
$ws.Range("L129").Value = 2606.21742
$ws.Range("M129").Value = 3934.625
$ws.Range("N129").Value = -12606.21742

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 30307906
$ws.Range("I132").Value = 35719050
$ws.Range("J132").Value = 5510
$ws.Range("K132").Value = 107157150
$ws.Range("L132").Value = 16530
$ws.Range("M132").Value = -107154620
$ws.Range("N132").Value = -21590

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5748.375
$ws.Range("I32").Value = 4368.4634
$ws.Range("J32").Value = 9520.134
$ws.Range("K32").Value = 4368.4634
$ws.Range("L32").Value = 9520.134
$ws.Range("M32").Value = -4081.4634
$ws.Range("N32").Value = -10094.134

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1957.9
$ws.Range("I132").Value = 1231.5217
$ws.Range("J132").Value = 4344.5713
$ws.Range("K132").Value = 3694.5651
$ws.Range("L132").Value = 13033.7139
$ws.Range("M132").Value = -1164.5651
$ws.Range("N132").Value = -18093.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 18999.5
$ws.Range("J14").Value = 18999.5
$ws.Range("L14").Value = 18999.5
$ws.Range("N14").Value = -19343.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 16844
$ws.Range("J17").Value = 16844
$ws.Range("L17").Value = 16844
$ws.Range("N17").Value = -17188

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 13199.454
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 13199.454
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 13199.454
$ws.Range("M15").Value = $null
$ws.Range("N15").Value = -13539.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 55556110
$ws.Range("I16").Value = 55556110
$ws.Range("K16").Value = 55556110
$ws.Range("M16").Value = -55555823

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2690.9375
$ws.Range("I31").Value = 1016.5
$ws.Range("J31").Value = 7714.25
$ws.Range("K31").Value = 1016.5
$ws.Range("L31").Value = 7714.25
$ws.Range("M31").Value = -721.5
$ws.Range("N31").Value = -8304.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2690.9375
$ws.Range("I34").Value = 1016.5
$ws.Range("J34").Value = 7714.25
$ws.Range("K34").Value = 1016.5
$ws.Range("L34").Value = 7714.25
$ws.Range("M34").Value = -814.5
$ws.Range("N34").Value = -8118.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3926
$ws.Range("I58").Value = 990.4
$ws.Range("J58").Value = 7071.2856
$ws.Range("K58").Value = 990.4
$ws.Range("L58").Value = 7071.2856
$ws.Range("M58").Value = -787.4
$ws.Range("N58").Value = -7477.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 55556110
$ws.Range("I113").Value = 55556110
$ws.Range("K113").Value = 55556110
$ws.Range("M113").Value = -55553940

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4114.926
$ws.Range("I132").Value = 3528.8096
$ws.Range("J132").Value = 6166.3335
$ws.Range("K132").Value = 10586.4288
$ws.Range("L132").Value = 18499.0005
$ws.Range("M132").Value = -8056.4288
$ws.Range("N132").Value = -23559.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3926
$ws.Range("I136").Value = 990.4
$ws.Range("J136").Value = 7071.2856
$ws.Range("K136").Value = 2971.2
$ws.Range("L136").Value = 21213.8568
$ws.Range("M136").Value = -421.1999999999998
$ws.Range("N136").Value = -26313.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 75499.75
$ws.Range("I4").Value = 300000
$ws.Range("J4").Value = 666.3333
$ws.Range("K4").Value = 900000
$ws.Range("L4").Value = 1998.9999
$ws.Range("M4").Value = -899888
$ws.Range("N4").Value = -2222.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 745.75
$ws.Range("I68").Value = 590
$ws.Range("J68").Value = 901.5
$ws.Range("K68").Value = 1770
$ws.Range("L68").Value = 2704.5
$ws.Range("M68").Value = -959
$ws.Range("N68").Value = -4326.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 745.75
$ws.Range("I71").Value = 590
$ws.Range("J71").Value = 901.5
$ws.Range("K71").Value = 5310
$ws.Range("L71").Value = 8113.5
$ws.Range("M71").Value = -1254
$ws.Range("N71").Value = -16225.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3125587
$ws.Range("I113").Value = 604.4286
$ws.Range("J113").Value = 6579515
$ws.Range("K113").Value = 1813.2858
$ws.Range("L113").Value = 19738545
$ws.Range("M113").Value = 356.7142000000001
$ws.Range("N113").Value = -19742885

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 3440
$ws.Range("I115").Value = 3000
$ws.Range("J115").Value = 3550
$ws.Range("K115").Value = 9000
$ws.Range("L115").Value = 10650
$ws.Range("M115").Value = -7825
$ws.Range("N115").Value = -13000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2358.6
$ws.Range("I130").Value = 1400
$ws.Range("J130").Value = 2997.6667
$ws.Range("K130").Value = 4200
$ws.Range("L130").Value = 8993.000100000001
$ws.Range("M130").Value = 820
$ws.Range("N130").Value = -19033.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 664.0599999999999
$ws.Range("I131").Value = 239.92308
$ws.Range("J131").Value = 813.08105
$ws.Range("K131").Value = 719.76924
$ws.Range("L131").Value = 2439.24315
$ws.Range("M131").Value = 4320.23076
$ws.Range("N131").Value = -12519.24315

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3543.5454
$ws.Range("I122").Value = 1577.75
$ws.Range("J122").Value = 8785.666999999999
$ws.Range("K122").Value = 4733.25
$ws.Range("L122").Value = 26357.001
$ws.Range("M122").Value = -2283.25
$ws.Range("N122").Value = -31257.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6199
$ws.Range("I7").Value = 5898.5
$ws.Range("K7").Value = 5898.5
$ws.Range("M7").Value = -5786.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6199
$ws.Range("I126").Value = 5898.5
$ws.Range("K126").Value = 17695.5
$ws.Range("M126").Value = -15225.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6231.375
$ws.Range("I132").Value = 2831.0667
$ws.Range("J132").Value = 11898.556
$ws.Range("K132").Value = 8493.2001
$ws.Range("L132").Value = 35695.66800000001
$ws.Range("M132").Value = -5963.2001
$ws.Range("N132").Value = -40755.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 762441
$ws.Range("I126").Value = 1877.4445
$ws.Range("J126").Value = 2131455.5
$ws.Range("K126").Value = 5632.333500000001
$ws.Range("L126").Value = 6394366.5
$ws.Range("M126").Value = -3162.333500000001
$ws.Range("N126").Value = -6399306.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5030.943
$ws.Range("I132").Value = 4824.4644
$ws.Range("K132").Value = 14473.3932
$ws.Range("M132").Value = -11943.3932

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14848.167
$ws.Range("I136").Value = 15596.143
$ws.Range("K136").Value = 46788.429
$ws.Range("M136").Value = -44238.429
Write-Host "Applied all cell updates"
